$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.383.99"
$ws.Range("E2").Value = "  +0.10%  "
$ws.Range("D3").Value = "1.911.18"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("D4").Formula = "=""1.007"""
$ws.Range("D4").Copy($ws.Range("D4"))
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = "  +0.72%  "
$ws.Range("D5").Formula = "=""324.83"""
$ws.Range("D5").Copy($ws.Range("D5"))
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +0.46%  "
$ws.Range("E6").Value = "  +0.67%  "
$ws.Range("D7").Formula = "=""0.4821"""
$ws.Range("D7").Copy($ws.Range("D7"))
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  +1.30%  "
$ws.Range("E8").Value = "  +0.41%  "
$ws.Range("D9").Formula = "=""0.08217"""
$ws.Range("D9").Copy($ws.Range("D9"))
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  +2.30%  "
$ws.Range("D10").Formula = "=""1.019"""
$ws.Range("D10").Copy($ws.Range("D10"))
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("D11").Formula = "=""23.45"""
$ws.Range("D11").Copy($ws.Range("D11"))
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").Value = "1.914.64"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").Formula = "=""6.049"""
$ws.Range("D13").Copy($ws.Range("D13"))
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  +2.01%  "
$ws.Range("D14").Formula = "=""7.206"""
$ws.Range("D14").Copy($ws.Range("D14"))
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  +2.18%  "
$ws.Range("D15").Formula = "=""91.00"""
$ws.Range("D15").Copy($ws.Range("D15"))
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("D16").Formula = "=""0.06801"""
$ws.Range("D16").Copy($ws.Range("D16"))
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  +1.86%  "
$ws.Range("E17").Value = "  +0.79%  "
$ws.Range("D18").Formula = "=""0.00001036"""
$ws.Range("D18").Copy($ws.Range("D18"))
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D19").Formula = "=""17.67"""
$ws.Range("D19").Copy($ws.Range("D19"))
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  +0.50%  "
$ws.Range("D20").Formula = "=""1.007"""
$ws.Range("D20").Copy($ws.Range("D20"))
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("D21").Value = "29.413.06"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").Formula = "=""5.626"""
$ws.Range("D22").Copy($ws.Range("D22"))
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  +1.86%  "
$ws.Range("E23").Value = "  +0.79%  "
$ws.Range("D24").Formula = "=""2.177"""
$ws.Range("D24").Copy($ws.Range("D24"))
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("D25").Value = "2.159.42"
$ws.Range("E25").Value = "  +0.61%  "
$ws.Range("D26").Formula = "=""6.588"""
$ws.Range("D26").Copy($ws.Range("D26"))
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("D27").Formula = "=""155.82"""
$ws.Range("D27").Copy($ws.Range("D27"))
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("D28").Formula = "=""19.98"""
$ws.Range("D28").Copy($ws.Range("D28"))
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  +1.11%  "
$ws.Range("D29").Formula = "=""2.103"""
$ws.Range("D29").Copy($ws.Range("D29"))
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  +0.48%  "
$ws.Range("D30").Formula = "=""120.25"""
$ws.Range("D30").Copy($ws.Range("D30"))
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  +2.09%  "
$ws.Range("D31").Formula = "=""1.019"""
$ws.Range("D31").Copy($ws.Range("D31"))
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("D32").Formula = "=""0.09554"""
$ws.Range("D32").Copy($ws.Range("D32"))
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  +1.36%  "
$ws.Range("D33").Formula = "=""5.599"""
$ws.Range("D33").Copy($ws.Range("D33"))
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  +4.62%  "
$ws.Range("D34").Formula = "=""3.547"""
$ws.Range("D34").Copy($ws.Range("D34"))
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("D35").Formula = "=""1.365"""
$ws.Range("D35").Copy($ws.Range("D35"))
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -1.08%  "
$ws.Range("D36").Formula = "=""0.02283"""
$ws.Range("D36").Copy($ws.Range("D36"))
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  +1.63%  "
$ws.Range("D37").Formula = "=""0.06102"""
$ws.Range("D37").Copy($ws.Range("D37"))
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("D38").Formula = "=""1.175"""
$ws.Range("D38").Copy($ws.Range("D38"))
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Formula = "=""0.5972"""
$ws.Range("D39").Copy($ws.Range("D39"))
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  +2.01%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Formula = "=""8.042"""
$ws.Range("D40").Copy($ws.Range("D40"))
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  +2.38%  "
$ws.Range("D41").Formula = "=""10.79"""
$ws.Range("D41").Copy($ws.Range("D41"))
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  +6.78%  "
$ws.Range("D42").Formula = "=""0.1845"""
$ws.Range("D42").Copy($ws.Range("D42"))
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  +0.42%  "
$ws.Range("D43").Formula = "=""2.409"""
$ws.Range("D43").Copy($ws.Range("D43"))
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  +1.33%  "
$ws.Range("D44").Formula = "=""1.280"""
$ws.Range("D44").Copy($ws.Range("D44"))
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -0.64%  "
$ws.Range("D45").Formula = "=""0.07610"""
$ws.Range("D45").Copy($ws.Range("D45"))
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -1.26%  "
$ws.Range("D46").Formula = "=""12.45"""
$ws.Range("D46").Copy($ws.Range("D46"))
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  +2.54%  "
$ws.Range("D47").Formula = "=""0.5568"""
$ws.Range("D47").Copy($ws.Range("D47"))
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  +1.25%  "
$ws.Range("D48").Formula = "=""1.951"""
$ws.Range("D48").Copy($ws.Range("D48"))
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  +1.55%  "
$ws.Range("D49").Formula = "=""117.60"""
$ws.Range("D49").Copy($ws.Range("D49"))
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  +3.87%  "
$ws.Range("D50").Formula = "=""2.418"""
$ws.Range("D50").Copy($ws.Range("D50"))
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  +3.79%  "
$ws.Range("D51").Formula = "=""72.10"""
$ws.Range("D51").Copy($ws.Range("D51"))
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  +0.52%  "
$excel.CutCopyMode = 0
